$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 19.49967933333333
$ws.Range("H2").Value = 58.499038
$ws.Range("I2").Value = 0.01453409262904611
$ws.Range("J2").Value = 0.01453409262904611
$ws.Range("M2").Value = 0.180428
$ws.Range("N2").Value = 0.541284
$ws.Range("O2").Value = 0.6724860231084607
$ws.Range("P2").Value = 0.6724860231084607
$ws.Range("Q2").Value = 3.518288142754667
$ws.Range("R2").Value = 31.664593284792
$ws.Range("S2").Value = 0.00977397415159721
$ws.Range("T2").Value = 0.00977397415159721

# Row 3
$ws.Range("G3").Value = 19.49967933333333
$ws.Range("H3").Value = 58.499038
$ws.Range("I3").Value = 0.01453409262904611
$ws.Range("J3").Value = 0.01453409262904611
$ws.Range("O3").Value = 0.06438936513852653
$ws.Range("P3").Value = 0.06438936513852653
$ws.Range("Q3").Value = 0.3368699602695555
$ws.Range("R3").Value = 3.031829642426
$ws.Range("S3").Value = 0.0009358409972488169
$ws.Range("T3").Value = 0.0009358409972488169

# Row 4
$ws.Range("G4").Value = 19.49967933333333
$ws.Range("H4").Value = 58.499038
$ws.Range("I4").Value = 0.01453409262904611
$ws.Range("J4").Value = 0.01453409262904611
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.07059633333333333
$ws.Range("N4").Value = 0.211789
$ws.Range("O4").Value = 0.2631246117530128
$ws.Range("P4").Value = 0.2631246117530128
$ws.Range("Q4").Value = 1.376605862109111
$ws.Range("R4").Value = 12.389452758982
$ws.Range("S4").Value = 0.003824277480200083
$ws.Range("T4").Value = 0.003824277480200083

# Row 5
$ws.Range("I5").Value = 0.2821439310161206
$ws.Range("J5").Value = 0.2821439310161206
$ws.Range("M5").Value = 0.180428
$ws.Range("N5").Value = 0.541284
$ws.Range("O5").Value = 0.6724860231084607
$ws.Range("P5").Value = 0.6724860231084607
$ws.Range("Q5").Value = 68.29897623332799
$ws.Range("R5").Value = 614.690786099952
$ws.Range("S5").Value = 0.1897378501132188
$ws.Range("T5").Value = 0.1897378501132188

# Row 6
$ws.Range("I6").Value = 0.2821439310161206
$ws.Range("J6").Value = 0.2821439310161206
$ws.Range("O6").Value = 0.06438936513852653
$ws.Range("P6").Value = 0.06438936513852653
$ws.Range("S6").Value = 0.01816706859581623
$ws.Range("T6").Value = 0.01816706859581623

# Row 7
$ws.Range("I7").Value = 0.2821439310161206
$ws.Range("J7").Value = 0.2821439310161206
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.07059633333333333
$ws.Range("N7").Value = 0.211789
$ws.Range("O7").Value = 0.2631246117530128
$ws.Range("P7").Value = 0.2631246117530128
$ws.Range("Q7").Value = 26.72344255045466
$ws.Range("R7").Value = 240.510982954092
$ws.Range("S7").Value = 0.07423901230708556
$ws.Range("T7").Value = 0.07423901230708557

# Row 8
$ws.Range("G8").Value = 481.5587156666667
$ws.Range("H8").Value = 1444.676147
$ws.Range("I8").Value = 0.3589299526510408
$ws.Range("J8").Value = 0.3589299526510408
$ws.Range("M8").Value = 0.180428
$ws.Range("N8").Value = 0.541284
$ws.Range("O8").Value = 0.6724860231084607
$ws.Range("P8").Value = 0.6724860231084607
$ws.Range("Q8").Value = 86.88667595030535
$ws.Range("R8").Value = 781.9800835527481
$ws.Range("S8").Value = 0.2413753764328065
$ws.Range("T8").Value = 0.2413753764328065

# Row 9
$ws.Range("G9").Value = 481.5587156666667
$ws.Range("H9").Value = 1444.676147
$ws.Range("I9").Value = 0.3589299526510408
$ws.Range("J9").Value = 0.3589299526510408
$ws.Range("O9").Value = 0.06438936513852653
$ws.Range("P9").Value = 0.06438936513852653
$ws.Range("Q9").Value = 8.319247852285445
$ws.Range("R9").Value = 74.873230670569
$ws.Range("S9").Value = 0.0231112717804019
$ws.Range("T9").Value = 0.0231112717804019

# Row 10
$ws.Range("G10").Value = 481.5587156666667
$ws.Range("H10").Value = 1444.676147
$ws.Range("I10").Value = 0.3589299526510408
$ws.Range("J10").Value = 0.3589299526510408
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.07059633333333333
$ws.Range("N10").Value = 0.211789
$ws.Range("O10").Value = 0.2631246117530128
$ws.Range("P10").Value = 0.2631246117530128
$ws.Range("Q10").Value = 33.99627961077589
$ws.Range("R10").Value = 305.966516496983
$ws.Range("S10").Value = 0.09444330443783239
$ws.Range("T10").Value = 0.09444330443783239

# Row 11
$ws.Range("G11").Value = 10.909999
$ws.Range("H11").Value = 32.729997
$ws.Range("I11").Value = 0.008131771468556478
$ws.Range("J11").Value = 0.008131771468556478
$ws.Range("M11").Value = 0.180428
$ws.Range("N11").Value = 0.541284
$ws.Range("O11").Value = 0.6724860231084607
$ws.Range("P11").Value = 0.6724860231084607
$ws.Range("Q11").Value = 1.968469299572
$ws.Range("R11").Value = 17.716223696148
$ws.Range("S11").Value = 0.005468502655716393
$ws.Range("T11").Value = 0.005468502655716393

# Row 12
$ws.Range("G12").Value = 10.909999
$ws.Range("H12").Value = 32.729997
$ws.Range("I12").Value = 0.008131771468556478
$ws.Range("J12").Value = 0.008131771468556478
$ws.Range("O12").Value = 0.06438936513852653
$ws.Range("P12").Value = 0.06438936513852653
$ws.Range("Q12").Value = 0.1884775060576666
$ws.Range("R12").Value = 1.696297554519
$ws.Range("S12").Value = 0.0005235996023119352
$ws.Range("T12").Value = 0.0005235996023119352

# Row 13
$ws.Range("G13").Value = 10.909999
$ws.Range("H13").Value = 32.729997
$ws.Range("I13").Value = 0.008131771468556478
$ws.Range("J13").Value = 0.008131771468556478
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.07059633333333333
$ws.Range("N13").Value = 0.211789
$ws.Range("O13").Value = 0.2631246117530128
$ws.Range("P13").Value = 0.2631246117530128
$ws.Range("Q13").Value = 0.7702059260703332
$ws.Range("R13").Value = 6.931853334633
$ws.Range("S13").Value = 0.00213966921052815
$ws.Range("T13").Value = 0.00213966921052815

# Row 14
$ws.Range("G14").Value = 98.48487833333333
$ws.Range("H14").Value = 295.454635
$ws.Range("I14").Value = 0.07340573759129181
$ws.Range("J14").Value = 0.07340573759129182
$ws.Range("M14").Value = 0.180428
$ws.Range("N14").Value = 0.541284
$ws.Range("O14").Value = 0.6724860231084607
$ws.Range("P14").Value = 0.6724860231084607
$ws.Range("Q14").Value = 17.76942962792667
$ws.Range("R14").Value = 159.92486665134
$ws.Range("S14").Value = 0.04936433254611106
$ws.Range("T14").Value = 0.04936433254611108

# Row 15
$ws.Range("G15").Value = 98.48487833333333
$ws.Range("H15").Value = 295.454635
$ws.Range("I15").Value = 0.07340573759129181
$ws.Range("J15").Value = 0.07340573759129182
$ws.Range("O15").Value = 0.06438936513852653
$ws.Range("P15").Value = 0.06438936513852653
$ws.Range("Q15").Value = 1.701391929793889
$ws.Range("R15").Value = 15.312527368145
$ws.Range("S15").Value = 0.004726548841028551
$ws.Range("T15").Value = 0.004726548841028552

# Row 16
$ws.Range("G16").Value = 98.48487833333333
$ws.Range("H16").Value = 295.454635
$ws.Range("I16").Value = 0.07340573759129181
$ws.Range("J16").Value = 0.07340573759129182
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.07059633333333333
$ws.Range("N16").Value = 0.211789
$ws.Range("O16").Value = 0.2631246117530128
$ws.Range("P16").Value = 0.2631246117530128
$ws.Range("Q16").Value = 6.952671299112777
$ws.Range("R16").Value = 62.574041692015
$ws.Range("S16").Value = 0.0193148562041522
$ws.Range("T16").Value = 0.0193148562041522

# Row 17
$ws.Range("G17").Value = 352.659012
$ws.Range("H17").Value = 1057.977036
$ws.Range("I17").Value = 0.2628545146439442
$ws.Range("J17").Value = 0.2628545146439442
$ws.Range("M17").Value = 0.180428
$ws.Range("N17").Value = 0.541284
$ws.Range("O17").Value = 0.6724860231084607
$ws.Range("P17").Value = 0.6724860231084607
$ws.Range("Q17").Value = 63.629560217136
$ws.Range("R17").Value = 572.6660419542239
$ws.Range("S17").Value = 0.1767659872090107
$ws.Range("T17").Value = 0.1767659872090107

# Row 18
$ws.Range("G18").Value = 352.659012
$ws.Range("H18").Value = 1057.977036
$ws.Range("I18").Value = 0.2628545146439442
$ws.Range("J18").Value = 0.2628545146439442
$ws.Range("O18").Value = 0.06438936513852653
$ws.Range("P18").Value = 0.06438936513852653
$ws.Range("Q18").Value = 6.092419538308
$ws.Range("R18").Value = 54.831775844772
$ws.Range("S18").Value = 0.01692503532171909
$ws.Range("T18").Value = 0.01692503532171909

# Row 19
$ws.Range("G19").Value = 352.659012
$ws.Range("H19").Value = 1057.977036
$ws.Range("I19").Value = 0.2628545146439442
$ws.Range("J19").Value = 0.2628545146439442
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.07059633333333333
$ws.Range("N19").Value = 0.211789
$ws.Range("O19").Value = 0.2631246117530128
$ws.Range("P19").Value = 0.2631246117530128
$ws.Range("Q19").Value = 24.896433164156
$ws.Range("R19").Value = 224.067898477404
$ws.Range("S19").Value = 0.06916349211321444
$ws.Range("T19").Value = 0.06916349211321444

